$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 308.25925
$ws.Range("I9").Value = 238.86363
$ws.Range("J9").Value = 613.6
$ws.Range("K9").Value = 238.86363
$ws.Range("L9").Value = 613.6
$ws.Range("M9").Value = -69.86363
$ws.Range("N9").Value = -951.6

$ws.Range("H38").Value = 858.8
$ws.Range("I38").Value = 64.666664
$ws.Range("K38").Value = 193.999992
$ws.Range("M38").Value = 178.000008

$ws.Range("H92").Value = 40000812
$ws.Range("I92").Value = 591.65
$ws.Range("K92").Value = 591.65
$ws.Range("M92").Value = 656.35

$ws.Range("H132").Value = 1620.9773
$ws.Range("I132").Value = 1562.7838
$ws.Range("K132").Value = 4688.3514
$ws.Range("M132").Value = -2158.3514

$ws.Range("H135").Value = 294869.03
$ws.Range("I135").Value = 323286.7
$ws.Range("K135").Value = 2909580.3
$ws.Range("M135").Value = -2907045.3

$ws.Range("H137").Value = 3572
$ws.Range("I137").Value = 3846.6
$ws.Range("K137").Value = 11539.8
$ws.Range("M137").Value = -8989.799999999999

$ws.Range("H138").Value = 2178607.8
$ws.Range("I138").Value = 959.85
$ws.Range("J138").Value = 3853721.5
$ws.Range("K138").Value = 2879.55
$ws.Range("L138").Value = 11561164.5
$ws.Range("M138").Value = 2260.45
$ws.Range("N138").Value = -11571444.5

$ws.Range("H141").Value = 1428.2727
$ws.Range("I141").Value = 865.7368
$ws.Range("K141").Value = 2597.2104
$ws.Range("M141").Value = 2582.7896

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H61").Value = 8582.799999999999
$ws.Range("J61").Value = 13973.637
$ws.Range("L61").Value = 13973.637
$ws.Range("N61").Value = -14397.637

$ws.Range("H64").Value = 33441
$ws.Range("J64").Value = 43000
$ws.Range("L64").Value = 43000
$ws.Range("N64").Value = -43496

$ws.Range("H67").Value = 33441
$ws.Range("J67").Value = 43000
$ws.Range("L67").Value = 43000
$ws.Range("N67").Value = -44716

$ws.Range("H74").Value = 25976.846
$ws.Range("I74").Value = 33778.844
$ws.Range("K74").Value = 33778.844
$ws.Range("M74").Value = -32904.844

$ws.Range("H77").Value = 25976.846
$ws.Range("I77").Value = 33778.844
$ws.Range("K77").Value = 168894.22
$ws.Range("M77").Value = -164526.22

$ws.Range("H97").Value = 3788700.8
$ws.Range("J97").Value = 41667816
$ws.Range("L97").Value = 41667816
$ws.Range("N97").Value = -41668808

$ws.Range("H102").Value = 1354
$ws.Range("I102").Value = 962.75
$ws.Range("J102").Value = 1980
$ws.Range("K102").Value = 962.75
$ws.Range("L102").Value = 1980
$ws.Range("M102").Value = 659.25
$ws.Range("N102").Value = -5224

$ws.Range("H132").Value = 4097.697
$ws.Range("I132").Value = 1977.1702
$ws.Range("K132").Value = 5931.5106
$ws.Range("M132").Value = -3401.5106

$ws.Range("H136").Value = 8582.799999999999
$ws.Range("J136").Value = 13973.637
$ws.Range("L136").Value = 41920.911
$ws.Range("N136").Value = -47020.911

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 36000
$ws.Range("J62").Value = 36000
$ws.Range("L62").Value = 36000
$ws.Range("N62").Value = -37372

$ws.Range("H65").Value = 36000
$ws.Range("J65").Value = 36000
$ws.Range("L65").Value = 108000
$ws.Range("N65").Value = -114864

$ws.Range("H105").Value = 4547
$ws.Range("I105").Value = 3994
$ws.Range("J105").Value = 4731.3335
$ws.Range("K105").Value = 3994
$ws.Range("L105").Value = 4731.3335
$ws.Range("M105").Value = -2247
$ws.Range("N105").Value = -8225.333500000001

$ws.Range("H134").Value = 5166.17
$ws.Range("I134").Value = 1885.2354
$ws.Range("K134").Value = 5655.706200000001
$ws.Range("M134").Value = -3120.706200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 29037
$ws.Range("J45").Value = 29037
$ws.Range("L45").Value = 29037
$ws.Range("N45").Value = -30223

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H112").Value = 69933
$ws.Range("J112").Value = 69933
$ws.Range("L112").Value = 69933
$ws.Range("N112").Value = -72887

$ws.Range("I132").Value = 2210.8774
$ws.Range("K132").Value = 6632.6322
$ws.Range("M132").Value = -4102.6322

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2098.639
$ws.Range("I131").Value = 689.8333
$ws.Range("K131").Value = 2069.4999
$ws.Range("M131").Value = 2970.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 29266.666
$ws.Range("J35").Value = 29266.666
$ws.Range("L35").Value = 29266.666
$ws.Range("N35").Value = -29862.666

$ws.Range("H58").Value = 65979.60000000001
$ws.Range("J58").Value = 65979.60000000001
$ws.Range("L58").Value = 65979.60000000001
$ws.Range("N58").Value = -66533.60000000001

$ws.Range("H63").Value = 48527
$ws.Range("J63").Value = 42069
$ws.Range("L63").Value = 42069
$ws.Range("N63").Value = -43441

$ws.Range("H66").Value = 48527
$ws.Range("J66").Value = 42069
$ws.Range("L66").Value = 126207
$ws.Range("N66").Value = -133071

$ws.Range("H70").Value = 71438680
$ws.Range("J70").Value = 11249.25
$ws.Range("L70").Value = 11249.25
$ws.Range("N70").Value = -11789.25

$ws.Range("H73").Value = 71438680
$ws.Range("J73").Value = 11249.25
$ws.Range("L73").Value = 11249.25
$ws.Range("N73").Value = -13121.25

$ws.Range("H80").Value = 94293.73
$ws.Range("I80").Value = 3533
$ws.Range("J80").Value = 253125
$ws.Range("K80").Value = 3533
$ws.Range("L80").Value = 253125
$ws.Range("M80").Value = -2535
$ws.Range("N80").Value = -255121

$ws.Range("H82").Value = 42500
$ws.Range("J82").Value = 42500
$ws.Range("L82").Value = 42500
$ws.Range("N82").Value = -43266

$ws.Range("H83").Value = 94293.73
$ws.Range("I83").Value = 3533
$ws.Range("J83").Value = 253125
$ws.Range("K83").Value = 17665
$ws.Range("L83").Value = 1265625
$ws.Range("M83").Value = -12673
$ws.Range("N83").Value = -1275609

$ws.Range("H85").Value = 42500
$ws.Range("J85").Value = 42500
$ws.Range("L85").Value = 42500
$ws.Range("N85").Value = -45152

$ws.Range("H132").Value = 5567.7715
$ws.Range("I132").Value = 2222.5715
$ws.Range("J132").Value = 10585.571
$ws.Range("K132").Value = 6667.7145
$ws.Range("L132").Value = 31756.713
$ws.Range("M132").Value = -4137.7145
$ws.Range("N132").Value = -36816.713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4646.758
$ws.Range("I93").Value = 3409.2273
$ws.Range("J93").Value = 7121.8184
$ws.Range("K93").Value = 3409.2273
$ws.Range("L93").Value = 7121.8184
$ws.Range("M93").Value = -2161.2273
$ws.Range("N93").Value = -9617.8184

$ws.Range("H100").Value = 4886.0835
$ws.Range("I100").Value = 4049.2
$ws.Range("J100").Value = 5483.857
$ws.Range("K100").Value = 4049.2
$ws.Range("L100").Value = 5483.857
$ws.Range("M100").Value = -3508.2
$ws.Range("N100").Value = -6565.857

$ws.Range("H139").Value = 93715
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 93715
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 93715
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -103995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 43998.6
$ws.Range("I15").Value = 34999.5
$ws.Range("K15").Value = 34999.5
$ws.Range("M15").Value = -34711.5

$ws.Range("H18").Value = 25529
$ws.Range("J18").Value = 25529
$ws.Range("L18").Value = 25529
$ws.Range("N18").Value = -25875

$ws.Range("H26").Value = 14
$ws.Range("J26").Value = 14
$ws.Range("L26").Value = 14
$ws.Range("N26").Value = -600

$ws.Range("H58").Value = 27994.5
$ws.Range("I58").Value = 27994.5
$ws.Range("K58").Value = 27994.5
$ws.Range("M58").Value = -27686.5

$ws.Range("H129").Value = 89428.5
$ws.Range("J129").Value = 89428.5
$ws.Range("L129").Value = 89428.5
$ws.Range("N129").Value = -99428.5

$ws.Range("H132").Value = 12204691
$ws.Range("I132").Value = 14709263
$ws.Range("K132").Value = 44127789
$ws.Range("M132").Value = -44125259

$ws.Range("H139").Value = 94823.125
$ws.Range("I139").Value = 90000
$ws.Range("J139").Value = 95512.14
$ws.Range("K139").Value = 90000
$ws.Range("L139").Value = 95512.14
$ws.Range("M139").Value = -84860
$ws.Range("N139").Value = -105792.14

$ws.Range("H141").Value = 88888.25
$ws.Range("J141").Value = 88888.25
$ws.Range("L141").Value = 88888.25
$ws.Range("N141").Value = -99248.25

Write-Host "Applied all 50 hunks"
